# Add a new column E ("fraction"/rate column) to the measles data sheet.
# E1 is a header cell (matches the style of the other header cells B1:D1),
# E2:E12 hold the new numeric rate values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (bold/centered/bordered header style) from D1 to E1
# so the new header cell matches the look of the existing header row.
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Header value for the new column.
$ws.Range("E1").Value = 3

# New data values for rows 2-12.
$ws.Range("E2").Value = 0.000413
$ws.Range("E3").Value = 0.002545
$ws.Range("E4").Value = 0.00228
$ws.Range("E5").Value = 0.000651
$ws.Range("E6").Value = 0.000169
$ws.Range("E7").Value = 0.000142
$ws.Range("E8").Value = 0.000111
$ws.Range("E9").Value = 0.000128
$ws.Range("E10").Value = 0.000109
$ws.Range("E11").Value = 0.000084
$ws.Range("E12").Value = 0.000107

Write-Output "Column E added successfully"
